$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.974.17"
$ws.Range("D3").Value = "1.640.02"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.66"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.53"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0882"
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "1.872.13"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.650.54"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.10"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.82"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "27.960.74"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.55"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.48"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.18"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.67"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.111"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "1.414.29"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.895"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +7.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.85"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.51"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "1.781.30"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.73"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("E51").Value = "  -1.65%  "
